$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SAP numbers (A2:A4) and the contenedores count (B2)
# to add storage type options for MP and MP1
$ws.Range("A2").Value = "1000009457A0"
$ws.Range("B2").Value = 3
$ws.Range("A3").Value = "1000011685A0"
$ws.Range("A4").Value = "1000013744A0"

# Move the active selection to C12
$ws.Range("C12").Select()
